# Quarterly rolling-window update for the Sergipe employment-sector series (g13.7).
# Each row's "Trimestre" (column C) and "Valor" (column D) shift to the NEXT
# quarter in the series; the newest quarter (01/07/2024) is appended at the
# bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new Trimestre text, new Valor number)
$updates = @(
    @(2, "01/07/2019", 25.6),
    @(3, "01/07/2019", 23.8),
    @(4, "01/07/2019", 13.4),
    @(5, "01/07/2019", 29.8),
    @(6, "01/07/2019", 3.8),
    @(7, "01/07/2019", 3.5),
    @(8, "01/10/2019", 27.4),
    @(9, "01/10/2019", 24.3),
    @(10, "01/10/2019", 13.8),
    @(11, "01/10/2019", 27.4),
    @(12, "01/10/2019", 4.2),
    @(13, "01/10/2019", 2.9),
    @(14, "01/01/2020", 28.8),
    @(15, "01/01/2020", 23.3),
    @(16, "01/01/2020", 13.2),
    @(17, "01/01/2020", 27.7),
    @(18, "01/01/2020", 4.6),
    @(19, "01/01/2020", 2.5),
    @(20, "01/04/2020", 28.2),
    @(21, "01/04/2020", 21),
    @(22, "01/04/2020", 15.9),
    @(23, "01/04/2020", 28.1),
    @(24, "01/04/2020", 4.6),
    @(25, "01/04/2020", 2.1),
    @(26, "01/07/2020", 28.6),
    @(27, "01/07/2020", 21.1),
    @(28, "01/07/2020", 15.2),
    @(29, "01/07/2020", 27.9),
    @(30, "01/07/2020", 4.1),
    @(31, "01/07/2020", 3.1),
    @(32, "01/10/2020", 25.2),
    @(33, "01/10/2020", 22.1),
    @(34, "01/10/2020", 14.3),
    @(35, "01/10/2020", 30),
    @(36, "01/10/2020", 4.9),
    @(37, "01/10/2020", 3.5),
    @(38, "01/01/2021", 27),
    @(39, "01/01/2021", 22.6),
    @(40, "01/01/2021", 14.9),
    @(41, "01/01/2021", 29.1),
    @(42, "01/01/2021", 3.1),
    @(43, "01/01/2021", 3.2),
    @(44, "01/04/2021", 25.3),
    @(45, "01/04/2021", 22.7),
    @(46, "01/04/2021", 16.9),
    @(47, "01/04/2021", 28.6),
    @(48, "01/04/2021", 3.2),
    @(49, "01/04/2021", 3.3),
    @(50, "01/07/2021", 23.7),
    @(51, "01/07/2021", 24.7),
    @(52, "01/07/2021", 16),
    @(53, "01/07/2021", 28.3),
    @(54, "01/07/2021", 4.4),
    @(55, "01/07/2021", 3),
    @(56, "01/10/2021", 25),
    @(57, "01/10/2021", 22.6),
    @(58, "01/10/2021", 14.5),
    @(59, "01/10/2021", 29.5),
    @(60, "01/10/2021", 4.8),
    @(61, "01/10/2021", 3.6),
    @(62, "01/01/2022", 26.4),
    @(63, "01/01/2022", 24.8),
    @(64, "01/01/2022", 13.7),
    @(65, "01/01/2022", 27.1),
    @(66, "01/01/2022", 5),
    @(67, "01/01/2022", 3),
    @(68, "01/04/2022", 27),
    @(69, "01/04/2022", 24.8),
    @(70, "01/04/2022", 13.4),
    @(71, "01/04/2022", 26),
    @(72, "01/04/2022", 5.8),
    @(73, "01/04/2022", 3),
    @(74, "01/07/2022", 27.3),
    @(75, "01/07/2022", 23.7),
    @(76, "01/07/2022", 15.2),
    @(77, "01/07/2022", 25.8),
    @(78, "01/07/2022", 4.9),
    @(79, "01/07/2022", 3.2),
    @(80, "01/10/2022", 28.7),
    @(81, "01/10/2022", 24.7),
    @(82, "01/10/2022", 15.5),
    @(83, "01/10/2022", 23.4),
    @(84, "01/10/2022", 4.3),
    @(85, "01/10/2022", 3.4),
    @(86, "01/01/2023", 28.6),
    @(87, "01/01/2023", 25.7),
    @(88, "01/01/2023", 15.7),
    @(89, "01/01/2023", 22.7),
    @(90, "01/01/2023", 4),
    @(91, "01/01/2023", 3.2),
    @(92, "01/04/2023", 28.2),
    @(93, "01/04/2023", 24.2),
    @(94, "01/04/2023", 15.6),
    @(95, "01/04/2023", 24.9),
    @(96, "01/04/2023", 4.1),
    @(97, "01/04/2023", 2.9),
    @(98, "01/07/2023", 27.1),
    @(99, "01/07/2023", 24.6),
    @(100, "01/07/2023", 15.4),
    @(101, "01/07/2023", 25.5),
    @(102, "01/07/2023", 4.2),
    @(103, "01/07/2023", 3.2),
    @(104, "01/10/2023", 28.2),
    @(105, "01/10/2023", 25.4),
    @(106, "01/10/2023", 15.7),
    @(107, "01/10/2023", 24.4),
    @(108, "01/10/2023", 3.5),
    @(109, "01/10/2023", 2.8),
    @(110, "01/01/2024", 29.3),
    @(111, "01/01/2024", 24.3),
    @(112, "01/01/2024", 15.2),
    @(113, "01/01/2024", 24.6),
    @(114, "01/01/2024", 3.6),
    @(115, "01/01/2024", 3),
    @(116, "01/04/2024", 30.2),
    @(117, "01/04/2024", 23.4),
    @(118, "01/04/2024", 15.3),
    @(119, "01/04/2024", 24.2),
    @(120, "01/04/2024", 4),
    @(121, "01/04/2024", 2.7),
    @(122, "01/07/2024", 30),
    @(123, "01/07/2024", 23.8),
    @(124, "01/07/2024", 14.7),
    @(125, "01/07/2024", 24.4),
    @(126, "01/07/2024", 4.5),
    @(127, "01/07/2024", 2.6)
)

foreach ($u in $updates) {
    $row = $u[0]
    $trimestre = $u[1]
    $valor = $u[2]

    # Force the "Trimestre" cell to stay literal text (matches the source file's
    # text cells) instead of letting Excel auto-convert the mm/dd/yyyy-looking
    # text into a real date serial number.
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $trimestre
    $cCell.ClearFormats()

    $ws.Cells.Item($row, 4).Value = $valor
}
